# Add two new "Title and Content" slides (layout 2) right after the
# existing title slide, matching the new slide2.xml / slide3.xml content.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2: "Background"
# ---------------------------------------------------------------------
$s2 = $p.Slides.Add(2, 2)

$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Background"

$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Text = "Described in Roy Fielding dissertation in 2000"
[void]$tr2.InsertAfter("`rRoots back to 1994")
[void]$tr2.InsertAfter("`rRoy was a key contributor to HTTP and URI")
[void]$tr2.InsertAfter("`rStyle of building APIs in distributed hypermedia systems")
[void]$tr2.InsertAfter("`rIt is not")
[void]$tr2.InsertAfter("`rParticular Framework or Implementation")
[void]$tr2.InsertAfter("`rSet of Standards")
[void]$tr2.InsertAfter("`r ")

$tr2.Paragraphs(2, 1).IndentLevel = 2
$tr2.Paragraphs(3, 1).IndentLevel = 2
$tr2.Paragraphs(6, 1).IndentLevel = 2
$tr2.Paragraphs(7, 1).IndentLevel = 2

$s2Last = $tr2.Paragraphs(8, 1)
$s2Last.Text = ""
$s2Last.IndentLevel = 2

# ---------------------------------------------------------------------
# Slide 3: "REST"
# ---------------------------------------------------------------------
$s3 = $p.Slides.Add(3, 2)

$s3.Shapes.Item(1).TextFrame.TextRange.Text = "REST"

$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$tr3.Text = "REprestational"
[void]$tr3.InsertAfter(" State Transfer")
[void]$tr3.InsertAfter("`rTransfers representations of resources in a particular State")
[void]$tr3.InsertAfter("`rSet of Constraints")
[void]$tr3.InsertAfter("`rBased on architectural style of WWW")
[void]$tr3.InsertAfter("`r ")
[void]$tr3.InsertAfter("`r ")

$tr3.Paragraphs(2, 1).IndentLevel = 2

$s3Empty1 = $tr3.Paragraphs(5, 1)
$s3Empty1.Text = ""

$s3Empty2 = $tr3.Paragraphs(6, 1)
$s3Empty2.Text = ""
$s3Empty2.IndentLevel = 2

Write-Host "Slides:" $p.Slides.Count
